# Update the cryptocurrency price/volume table (columns D and E) on Sheet1
# to reflect the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are stored as plain text (they use "."
# as a thousands separator, e.g. "28.199.30", and keep trailing zeros,
# e.g. "1.0000"), so force a Text format before assigning each value and
# then restore the cell's original ("Normal") style so no formatting changes leak in.
$prices = [ordered]@{
    "D2" = "28.199.30"
    "D3" = "1.810.27"
    "D5" = "339.17"
    "D6" = "0.9996"
    "D7" = "0.3936"
    "D8" = "0.3503"
    "D9" = "48.11"
    "D10" = "1.175"
    "D11" = "0.07550"
    "D12" = "1.0000"
    "D13" = "22.13"
    "D14" = "6.515"
    "D15" = "1.813.01"
    "D16" = "7.166"
    "D17" = "0.00001103"
    "D18" = "0.06704"
    "D19" = "85.28"
    "D20" = "0.9999"
    "D22" = "6.565"
    "D23" = "28.187.58"
    "D24" = "12.39"
    "D26" = "21.53"
    "D27" = "1.480"
    "D29" = "154.27"
    "D30" = "2.015.08"
    "D31" = "136.12"
    "D32" = "6.228"
    "D33" = "4.012"
    "D34" = "0.08860"
    "D35" = "13.25"
    "D36" = "0.02438"
    "D37" = "0.6934"
    "D38" = "5.463"
    "D39" = "0.06524"
    "D40" = "1.606"
    "D42" = "1.261"
    "D44" = "14.64"
    "D45" = "0.6431"
    "D46" = "0.9992"
    "D47" = "3.873"
    "D48" = "2.150"
    "D49" = "131.22"
    "D50" = "0.07189"
    "D51" = "80.30"
}
foreach ($addr in $prices.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $prices[$addr]
    $cell.Style = "Normal"
}

# Column E ("Volume(1h)") values are plain padded percentage strings.
$volumes = [ordered]@{
    "E2" = "  +3.49%  "
    "E3" = "  +1.39%  "
    "E4" = "  -0.49%  "
    "E5" = "  +0.83%  "
    "E6" = "  -0.62%  "
    "E7" = "  +3.80%  "
    "E8" = "  +2.13%  "
    "E9" = "  -0.67%  "
    "E10" = "  -1.42%  "
    "E11" = "  +1.25%  "
    "E12" = "  -0.53%  "
    "E13" = "  +1.57%  "
    "E14" = "  +1.26%  "
    "E15" = "  +1.64%  "
    "E16" = "  +1.43%  "
    "E17" = "  +0.49%  "
    "E18" = "  +0.33%  "
    "E19" = "  +0.95%  "
    "E20" = "  -0.27%  "
    "E21" = "  +2.40%  "
    "E22" = "  +0.64%  "
    "E23" = "  +3.56%  "
    "E24" = "  -0.31%  "
    "E25" = "  -1.19%  "
    "E26" = "  +0.92%  "
    "E27" = "  -0.96%  "
    "E28" = "  -0.56%  "
    "E29" = "  +0.84%  "
    "E30" = "  +1.38%  "
    "E31" = "  +2.45%  "
    "E32" = "  +3.12%  "
    "E33" = "  -1.26%  "
    "E34" = "  +2.50%  "
    "E35" = "  +1.48%  "
    "E36" = "  +4.37%  "
    "E37" = "  +1.23%  "
    "E38" = "  +0.23%  "
    "E39" = "  +2.15%  "
    "E40" = "  -2.74%  "
    "E41" = "  +1.09%  "
    "E42" = "  +0.04%  "
    "E43" = "  -2.93%  "
    "E44" = "  +1.06%  "
    "E45" = "  +0.37%  "
    "E46" = "  -0.55%  "
    "E47" = "  +0.40%  "
    "E48" = "  +1.09%  "
    "E49" = "  +1.79%  "
    "E50" = "  +0.22%  "
    "E51" = "  +1.51%  "
}
foreach ($addr in $volumes.Keys) {
    $ws.Range($addr).Value = $volumes[$addr]
}

Write-Host "Updated cryptos list"
